$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date in column C for rows 2-8 from 2023-10-08 to 2023-10-09
foreach ($r in 2..8) {
    $ws.Cells.Item($r, 3).Value = 45208
}
